# Applies the "kreise.xlsx" update described in the commit message:
#  - Nullpunkt Kreis Hamm aktualisiert: new NP X / NP Y values and updated NP Hinweis text
#  - Vermessungsnetz Camen (Dortmund row F33) ergänzt
#  - active cell selection moved to G10
#  - minor formatting cleanups (iterateDelta / window position are metadata,
#    not reachable through the object model and are left to Excel defaults)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kreise")

# --- Update row 9 (Hamm): NP X, NP Y, and NP Hinweis -----------------------
$ws.Range("D9").Value = 409239.2
$ws.Range("E9").Value = 5710123.6100000003
$ws.Range("G9").Value = "Nullpunkt wahrscheinlich Stadtkiche Unna, für einige Bgmstr. bestätigt. Rotation berechnet über Bgmstr. Unna"

# --- Update row 33 (Wiedenbrück / Vermessungsnetz Camen): Net Rotation -----
$ws.Range("F33").Value = 1.575
# Pick up the same cell formatting used by the other filled-in "Net Rotation"
# cells in the table (e.g. F9) instead of the default blank-row border.
$ws.Range("F9").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the active cell selection --------------------------------------
$ws.Range("G10").Select()
